# Update countries & provincias Spain
# This script applies the latest COVID-19 data refresh:
#  - updates the "last updated" timestamp
#  - updates case totals for several countries
#  - a few countries swapped rank/row position as their totals crossed over

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 18 de Septiembre de 2020 a las 20:29"

$rowsData = @(
    @{ Row = 4; A = "Estados Unidos"; B = 6895537; C = 20941; D = 4171374; E = 2521495; F = 0; G = 455; H = 202668 },
    @{ Row = 5; A = "India"; B = 5305387; C = 92701; D = 4203017; E = 1016749; F = 0; G = 1217; H = 85621 },
    @{ Row = 6; A = "Brasil"; B = 4466828; C = 9385; D = 3753082; E = 578543; F = 0; G = 172; H = 135203 },
    @{ Row = 11; A = "España"; B = 659334; C = 4697; D = 0; E = 0; F = 0; G = 90; H = 30495 },
    @{ Row = 12; A = "Sudafrica"; B = 655572; C = 0; D = 585303; E = 54497; F = 0; G = 0; H = 15772 },
    @{ Row = 15; A = "Francia"; B = 428696; C = 13215; D = 91574; E = 305873; F = 0; G = 154; H = 31249 },
    @{ Row = 16; A = "Iran"; B = 416198; C = 3049; D = 355505; E = 36741; F = 0; G = 144; H = 23952 },
    @{ Row = 29; A = "Canada"; B = 141566; C = 699; D = 123512; E = 8853; F = 0; G = 1; H = 9201 },
    @{ Row = 40; A = "Marruecos"; B = 97264; C = 2760; D = 76690; E = 18819; F = 0; G = 41; H = 1755 },
    @{ Row = 51; A = "Etiopia"; B = 67515; C = 602; D = 27638; E = 38805; F = 0; G = 12; H = 1072 },
    @{ Row = 52; A = "Portugal"; B = 67176; C = 780; D = 45053; E = 20229; F = 0; G = 6; H = 1894 },
    @{ Row = 54; A = "Barein"; B = 63189; C = 0; D = 56087; E = 6882; F = 0; G = 3; H = 220 },
    @{ Row = 60; A = "Argelia"; B = 49413; C = 219; D = 34818; E = 12936; F = 0; G = 5; H = 1659 },
    @{ Row = 61; A = "Suiza"; B = 49283; C = 488; D = 39900; E = 7338; F = 0; G = 3; H = 2045 },
    @{ Row = 75; A = "Libano"; B = 27518; C = 750; D = 10739; E = 16498; F = 0; G = 18; H = 281 },
    @{ Row = 76; A = "El Salvador"; B = 27346; C = 97; D = 20825; E = 5717; F = 0; G = 3; H = 804 },
    @{ Row = 77; A = "Australia"; B = 26861; C = 48; D = 23855; E = 2169; F = 0; G = 5; H = 837 },
    @{ Row = 93; A = "Noruega"; B = 12708; C = 64; D = 10371; E = 2070; F = 0; G = 1; H = 267 },
    @{ Row = 95; A = "Consejo Danes para los Refugiados"; B = 10456; C = 14; D = 9863; E = 325; F = 0; G = 1; H = 268 },
    @{ Row = 139; A = "Sri Lanka"; B = 3281; C = 5; D = 3060; E = 208; F = 0; G = 0; H = 13 },
    @{ Row = 143; A = "Mali"; B = 2991; C = 25; D = 2332; E = 531; F = 0; G = 0; H = 128 },
    @{ Row = 153; A = "Yemen"; B = 2024; C = 2; D = 1221; E = 218; F = 0; G = 0; H = 585 },
    @{ Row = 214; A = "Islas Malvinas"; B = 13; C = 0; D = 13; E = 0; F = 0; G = 0; H = 0 },
    @{ Row = 215; A = "Montserrat"; B = 13; C = 0; D = 12; E = 0; F = 0; G = 0; H = 1 }
)

foreach ($r in $rowsData) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}
